$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = "Última actualización: 04:28:33"
$ws.Cells.Item(3, 1).Value = "Total filas: 16"
$ws.Cells.Item(6, 1).Value = "03:58:57"
$ws.Cells.Item(6, 2).Value = "04:01"
$ws.Cells.Item(6, 3).Value = "81_EL PELIGRO"
$ws.Cells.Item(6, 4).Value = 3
$ws.Cells.Item(6, 5).Value = "LP1912"
$ws.Cells.Item(7, 1).Value = "03:58:57"
$ws.Cells.Item(7, 2).Value = "04:45"
$ws.Cells.Item(7, 3).Value = "215A_EL PATO"
$ws.Cells.Item(7, 4).Value = 47
$ws.Cells.Item(7, 5).Value = "LP1912"
$ws.Cells.Item(8, 1).Value = "04:28:33"
$ws.Cells.Item(8, 2).Value = "04:46"
$ws.Cells.Item(8, 3).Value = "215A_EL PATO"
$ws.Cells.Item(8, 4).Value = 18
$ws.Cells.Item(8, 5).Value = "LP1912"
$ws.Cells.Item(9, 1).Value = "03:58:57"
$ws.Cells.Item(9, 2).Value = "04:53"
$ws.Cells.Item(9, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(9, 4).Value = 55
$ws.Cells.Item(9, 5).Value = "LP1912"
$ws.Cells.Item(10, 1).Value = "03:58:57"
$ws.Cells.Item(10, 2).Value = "05:16"
$ws.Cells.Item(10, 3).Value = "17_ROMERO"
$ws.Cells.Item(10, 4).Value = 78
$ws.Cells.Item(10, 5).Value = "LP1912"
$ws.Cells.Item(11, 1).Value = "03:58:57"
$ws.Cells.Item(11, 2).Value = "05:21"
$ws.Cells.Item(11, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(11, 4).Value = 83
$ws.Cells.Item(11, 5).Value = "LP1912"
$ws.Cells.Item(12, 1).Value = "04:28:33"
$ws.Cells.Item(12, 2).Value = "05:22"
$ws.Cells.Item(12, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(12, 4).Value = 54
$ws.Cells.Item(12, 5).Value = "LP1912"
$ws.Cells.Item(13, 1).Value = "03:58:57"
$ws.Cells.Item(13, 2).Value = "05:34"
$ws.Cells.Item(13, 3).Value = "215B_EL PATO"
$ws.Cells.Item(13, 4).Value = 96
$ws.Cells.Item(13, 5).Value = "LP1912"
$ws.Cells.Item(14, 1).Value = "03:58:57"
$ws.Cells.Item(14, 2).Value = "05:46"
$ws.Cells.Item(14, 3).Value = "15_ABASTO"
$ws.Cells.Item(14, 4).Value = 108
$ws.Cells.Item(14, 5).Value = "LP1912"
$ws.Cells.Item(15, 1).Value = "03:58:57"
$ws.Cells.Item(15, 2).Value = "05:53"
$ws.Cells.Item(15, 3).Value = "10_OLMOS"
$ws.Cells.Item(15, 4).Value = 115
$ws.Cells.Item(15, 5).Value = "LP1912"
$ws.Cells.Item(16, 1).Value = "04:28:33"
$ws.Cells.Item(16, 2).Value = "05:54"
$ws.Cells.Item(16, 3).Value = "10_OLMOS"
$ws.Cells.Item(16, 4).Value = 86
$ws.Cells.Item(16, 5).Value = "LP1912"
$ws.Cells.Item(17, 1).Value = "04:28:33"
$ws.Cells.Item(17, 2).Value = "06:04"
$ws.Cells.Item(17, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(17, 4).Value = 96
$ws.Cells.Item(17, 5).Value = "LP1912"
$ws.Cells.Item(18, 1).Value = "04:28:33"
$ws.Cells.Item(18, 2).Value = "06:11"
$ws.Cells.Item(18, 3).Value = "215A_EL PATO"
$ws.Cells.Item(18, 4).Value = 103
$ws.Cells.Item(18, 5).Value = "LP1912"
$ws.Cells.Item(19, 1).Value = "04:28:33"
$ws.Cells.Item(19, 2).Value = "06:14"
$ws.Cells.Item(19, 3).Value = "225_HARAS DEL SUR"
$ws.Cells.Item(19, 4).Value = 106
$ws.Cells.Item(19, 5).Value = "LP1912"
$ws.Cells.Item(20, 1).Value = "04:28:33"
$ws.Cells.Item(20, 2).Value = "06:21"
$ws.Cells.Item(20, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(20, 4).Value = 113
$ws.Cells.Item(20, 5).Value = "LP1912"
$ws.Cells.Item(21, 1).Value = "04:28:33"
$ws.Cells.Item(21, 2).Value = "06:27"
$ws.Cells.Item(21, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(21, 4).Value = 119
$ws.Cells.Item(21, 5).Value = "LP1912"

$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = "Última actualización: 04:28:33"
$ws.Cells.Item(3, 1).Value = "Total filas: 4"
$ws.Cells.Item(6, 1).Value = "03:58:57"
$ws.Cells.Item(6, 2).Value = "04:45"
$ws.Cells.Item(6, 3).Value = "215A_EL PATO"
$ws.Cells.Item(6, 4).Value = 47
$ws.Cells.Item(6, 5).Value = "LP1912"
$ws.Cells.Item(7, 1).Value = "04:28:33"
$ws.Cells.Item(7, 2).Value = "04:46"
$ws.Cells.Item(7, 3).Value = "215A_EL PATO"
$ws.Cells.Item(7, 4).Value = 18
$ws.Cells.Item(7, 5).Value = "LP1912"
$ws.Cells.Item(8, 1).Value = "03:58:57"
$ws.Cells.Item(8, 2).Value = "05:34"
$ws.Cells.Item(8, 3).Value = "215B_EL PATO"
$ws.Cells.Item(8, 4).Value = 96
$ws.Cells.Item(8, 5).Value = "LP1912"
$ws.Cells.Item(9, 1).Value = "04:28:33"
$ws.Cells.Item(9, 2).Value = "06:11"
$ws.Cells.Item(9, 3).Value = "215A_EL PATO"
$ws.Cells.Item(9, 4).Value = 103
$ws.Cells.Item(9, 5).Value = "LP1912"

$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = "Última actualización: 04:28:33"
$ws.Cells.Item(3, 1).Value = "Total filas: 3"
$ws.Cells.Item(6, 1).Value = "03:58:57"
$ws.Cells.Item(6, 2).Value = "05:43"
$ws.Cells.Item(6, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(6, 4).Value = 105
$ws.Cells.Item(6, 5).Value = "L6173"
$ws.Cells.Item(7, 1).Value = "04:28:33"
$ws.Cells.Item(7, 2).Value = "05:44"
$ws.Cells.Item(7, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(7, 4).Value = 76
$ws.Cells.Item(7, 5).Value = "L6173"
$ws.Cells.Item(8, 1).Value = "04:28:33"
$ws.Cells.Item(8, 2).Value = "06:09"
$ws.Cells.Item(8, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(8, 4).Value = 101
$ws.Cells.Item(8, 5).Value = "L6173"
